$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "hello",
    "what is the order status",
    "what is the roode rstatsus",
    "what is order status in july",
    "what is order status in july",
    "what is order number of the latest order i plaved in july",
    "wht is its status",
    "what is the status of  1823383.",
    "list all the orders i placed in july"
)

$startRow = 187
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("A$row").Value = $values[$i]
}
